$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First "fridlysta arter" comment (Ur FSC-standarden section):
#    strip the italic formatting and add a trailing period.
# ------------------------------------------------------------------
$oldText1 = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"
$newText1 = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = $newText1
    $rng1.Font.Italic = $false
}

# ------------------------------------------------------------------
# 2) 6.4 paragraph: drop the stray trailing space on the run text.
# ------------------------------------------------------------------
$oldText2 = "Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten. Det ska ske genom avsättningar, andra skyddade områden och genom att skapa konnektivitet och/eller genom andra direkta åtgärder som gynnar dessa arters överlevnad och livskraft. Åtgärderna ska stå i förhållande till brukandets skala, intensitet och risk, samt till sällsynta och hotade arters bevarandestatus och ekologiska krav. Certifikatsinnehavaren ska beakta den geografiska spridningen och ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas. "
$newText2 = "Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten. Det ska ske genom avsättningar, andra skyddade områden och genom att skapa konnektivitet och/eller genom andra direkta åtgärder som gynnar dessa arters överlevnad och livskraft. Åtgärderna ska stå i förhållande till brukandets skala, intensitet och risk, samt till sällsynta och hotade arters bevarandestatus och ekologiska krav. Certifikatsinnehavaren ska beakta den geografiska spridningen och ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas."

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Text = $newText2
}

# ------------------------------------------------------------------
# 3) Remove the two "6.4.1 Följande biotoper..." / "b) nyckelbiotoper..."
#    paragraphs, then renumber the remaining "6.4.1" clause to "6.4.3".
# ------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Contains("Följande biotoper undantas")) {
        $pFirst = $p
        $pSecond = $paras.Item($i + 1)
        $delRange = $d.Range($pFirst.Range.Start, $pSecond.Range.End)
        $delRange.Delete()
        break
    }
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("6.4.1 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Text = "6.4.3 "
}

# ------------------------------------------------------------------
# 4) Update the letter date in the first-page header.
# ------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-10-22", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-25", 2)
